$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new match rows appended at the bottom of the table (rows 78 and 79),
# following the same column layout as the existing data rows (2-77):
# Indice, pais, torneio, temporada, data_partida, home, home_ft_gols, away,
# away_ft_gols, home_opening_odds, home_opening_data_hora, home_closing_odds,
# home_closing_data_hora, draw_opening_odds, draw_opening_data_hora,
# draw_closing_odds, draw_closing_data_hora, away_opening_odds,
# away_opening_data_hora, away_closing_odds, away_closing_data_hora, url_partida

$rows = @(
    @{
        Row = 78
        Indice = 77
        DataPartida = 45254.5
        Home = "Buriram"
        HomeGols = 0
        Away = "Ratchaburi"
        AwayGols = 0
        HomeOpenOdds = 1.39
        HomeOpenData = "17/11/2023 12:12"
        HomeCloseOdds = 1.42
        HomeCloseData = "24/11/2023 11:56"
        DrawOpenOdds = 4.85
        DrawOpenData = "17/11/2023 12:12"
        DrawCloseOdds = 4.86
        DrawCloseData = "24/11/2023 11:59"
        AwayOpenOdds = 7.42
        AwayOpenData = "17/11/2023 12:12"
        AwayCloseOdds = 7.09
        AwayCloseData = "24/11/2023 11:59"
        Url = "https://www.betexplorer.com/football/thailand/thai-league-1/buriram-united-f-c-ratchaburi/jDH5hVpg/"
    },
    @{
        Row = 79
        Indice = 78
        DataPartida = 45254.54166666666
        Home = "Bangkok Utd"
        HomeGols = 1
        Away = "Nakhon Pathom"
        AwayGols = 1
        HomeOpenOdds = 1.22
        HomeOpenData = "17/11/2023 13:12"
        HomeCloseOdds = 1.21
        HomeCloseData = "24/11/2023 11:56"
        DrawOpenOdds = 6.73
        DrawOpenData = "17/11/2023 13:12"
        DrawCloseOdds = 6.38
        DrawCloseData = "24/11/2023 11:56"
        AwayOpenOdds = 10.86
        AwayOpenData = "17/11/2023 13:12"
        AwayCloseOdds = 13.38
        AwayCloseData = "24/11/2023 11:56"
        Url = "https://www.betexplorer.com/football/thailand/thai-league-1/bangkok-utd-nakhon-pathom/Islz8Apt/"
    }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Mirror the formatting (styles on column A / E) of the previous row.
    $ws.Range("A" + ($row - 1) + ":V" + ($row - 1)).Copy() | Out-Null
    $ws.Range("A" + $row + ":V" + $row).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $r.Indice
    $ws.Cells.Item($row, 2).Value = "thailand"
    $ws.Cells.Item($row, 3).Value = "thai-league-1"
    $ws.Cells.Item($row, 4).Value = "2023-2024"

    $ws.Cells.Item($row, 5).Value = $r.DataPartida

    $ws.Cells.Item($row, 6).Value = $r.Home
    $ws.Cells.Item($row, 7).Value = $r.HomeGols
    $ws.Cells.Item($row, 8).Value = $r.Away
    $ws.Cells.Item($row, 9).Value = $r.AwayGols

    $ws.Cells.Item($row, 10).Value = $r.HomeOpenOdds
    $ws.Cells.Item($row, 11).Value = $r.HomeOpenData
    $ws.Cells.Item($row, 12).Value = $r.HomeCloseOdds
    $ws.Cells.Item($row, 13).Value = $r.HomeCloseData

    $ws.Cells.Item($row, 14).Value = $r.DrawOpenOdds
    $ws.Cells.Item($row, 15).Value = $r.DrawOpenData
    $ws.Cells.Item($row, 16).Value = $r.DrawCloseOdds
    $ws.Cells.Item($row, 17).Value = $r.DrawCloseData

    $ws.Cells.Item($row, 18).Value = $r.AwayOpenOdds
    $ws.Cells.Item($row, 19).Value = $r.AwayOpenData
    $ws.Cells.Item($row, 20).Value = $r.AwayCloseOdds
    $ws.Cells.Item($row, 21).Value = $r.AwayCloseData

    $ws.Cells.Item($row, 22).Value = $r.Url
}

$excel.CutCopyMode = 0
